$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell -> new value (odds updates as per commit "Atualizando o arquivo XLSX")
$updates = @{
    "Q2"  = 2.1
    "R2"  = 1.73
    "U2"  = 1.8
    "V2"  = 1.95
    "AT2" = 2.63

    "Q3"  = 2.3
    "R3"  = 1.62
    "U3"  = 1.95
    "V3"  = 1.8
    "BD3" = 151

    "Q4"  = 2.3
    "R4"  = 1.62
    "V4"  = 1.62

    "J12" = 2.88

    "M17" = 1.03
    "O17" = 1.22

    "M18" = 1.01
    "O18" = 1.1

    "M19" = 1.03
    "O19" = 1.17

    "M20" = 1.03
    "O20" = 1.19

    "M21" = 1.01
    "O21" = 1.11
    "Q21" = 1.48

    "O22" = 1.07
    "Q22" = 1.33

    "M23" = 1.05
    "O23" = 1.37
    "R23" = 1.62

    "S25" = 1.4
    "S26" = 1.33
    "S27" = 1.5
    "S28" = 1.4
    "S30" = 1.4
    "S31" = 1.29
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
